$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create WO")

# Update the Item Number value in B2 (was "Pro-Lot Track (Lot Track)")
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Move the active selection to B2, matching the saved sheet view state
$ws.Range("B2").Select()
